# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$updates = @{
    7  = 560
    9  = 6815
    13 = 177
    16 = 16226
    17 = 1595
    22 = 11382
    24 = 1026
    26 = 322
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates4 = @{
    7  = 560
    10 = 6815
    14 = 177
    18 = 16226
    19 = 1595
    26 = 11382
    28 = 1026
    30 = 322
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}

$wb.Save()
